$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "Classified.v2 MPMLP" -> "Classified.v4 MPMLP" everywhere it appears (column B) ---
$ws.Range("B4").Value  = "Classified.v4 MPMLP"
$ws.Range("B7").Value  = "Classified.v4 MPMLP"
$ws.Range("B10").Value = "Classified.v4 MPMLP"
$ws.Range("B13").Value = "Classified.v4 MPMLP"
$ws.Range("B16").Value = "Classified.v4 MPMLP"

# --- Fill in previously empty MPMLP rows (2, 5, 8, 11, 14) with their metrics ---
# Row 2: bfs-10 / MPMLP
$ws.Range("C2").Value = "Overfitting(1)"
$ws.Range("D2").Value = 5.3627000000000002
$ws.Range("E2").Value = 0.37390000000000001

# Row 5: 471.omnetpp-s2 / MPMLP
$ws.Range("C5").Value = 0.5867
$ws.Range("D5").Value = 46.8187
$ws.Range("E5").Value = 0.61829999999999996

# Row 8: 482.sphinx3-s0 / MPMLP
$ws.Range("C8").Value = 0.69720000000000004
$ws.Range("D8").Value = 15.7012
$ws.Range("E8").Value = 0.58250000000000002

# Row 11: 605.mcf-s8 / MPMLP
$ws.Range("C11").Value = 0.13450000000000001
$ws.Range("D11").Value = 32.976399999999998
$ws.Range("E11").Value = 0.037600000000000001

# Row 14: 623.xalancbmk-s1 / MPMLP
$ws.Range("C14").Value = 0.95440000000000003
$ws.Range("D14").Value = 12.086
$ws.Range("E14").Value = 0.7177

# --- Cosmetic changes ---
# Column C width widened (target stored width ~28.33203125 characters)
$ws.Columns.Item(3).ColumnWidth = 27.5

# Selected cell moved to H14
$ws.Range("H14").Select()
